$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.722.66"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.733.17"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.52"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.28"
$ws.Range("E6").Value = "  +6.30%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.732.86"
$ws.Range("E9").Value = "  +3.13%  "
$ws.Range("E10").Value = "  +3.02%  "
$ws.Range("E11").Value = "  +4.81%  "
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.84"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.230.27"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.765.50"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.694.30"
$ws.Range("E18").Value = "  +2.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.94"
$ws.Range("E19").Value = "  +4.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "376.22"
$ws.Range("E20").Value = "  +4.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.68"
$ws.Range("E21").Value = "  +3.55%  "
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("E23").Value = "  +5.28%  "
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.68"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("E27").Value = "  +3.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.871.50"
$ws.Range("E28").Value = "  +3.04%  "
$ws.Range("E29").Value = "  +2.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "591.64"
$ws.Range("E30").Value = "  +5.14%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +3.79%  "
$ws.Range("E33").Value = "  +4.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.00"
$ws.Range("E35").Value = "  +3.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.62"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.22"
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.381"
$ws.Range("E40").Value = "  +2.96%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.91"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.50"
$ws.Range("E42").Value = "  +2.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.67"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.09"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.606"
$ws.Range("E48").Value = "  +5.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "155.76"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("E50").Value = "  +3.66%  "
$ws.Range("E51").Value = "  +6.06%  "

Write-Output "Applied cryptos update"
